$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '54.573.54'
$ws.Range('E2').Value = '  -6.65%  '

# Row 3
$ws.Range('D3').Value = '2.420.18'
$ws.Range('E3').Value = '  -10.90%  '

# Row 4
$ws.Range('E4').Value = '  +0.02%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '468.07'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -6.87%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '132.56'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -5.64%  '

# Row 7
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.997'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -0.07%  '

# Row 8
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.493'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -6.64%  '

# Row 9
$ws.Range('D9').Value = '2.434.74'
$ws.Range('E9').Value = '  -10.70%  '

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0954'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -8.74%  '

# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '5.33'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -11.96%  '

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.315'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -9.11%  '

# Row 13
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.122'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -3.89%  '

# Row 14
$ws.Range('D14').Value = '2.844.64'
$ws.Range('E14').Value = '  -10.96%  '

# Row 15
$ws.Range('D15').Value = '54.514.88'
$ws.Range('E15').Value = '  -6.91%  '

# Row 16
$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.0000133'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -1.17%  '

# Row 17
$ws.Range('B17').Value = 'Avalanche'
$ws.Range('C17').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '19.71'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -8.87%  '

# Row 18
$ws.Range('D18').Value = '2.440.57'
$ws.Range('E18').Value = '  -9.39%  '

# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '4.20'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -11.42%  '

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '311.02'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -8.78%  '

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '9.51'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -13.02%  '

# Row 22
$ws.Range('E22').Value = '  -0.14%  '

# Row 23
$ws.Range('E23').Value = '  -0.08%  '

# Row 24
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '5.39'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -13.69%  '

# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '56.36'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -10.28%  '

# Row 26
$ws.Range('E26').Value = '  +1.09%  '

# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.386'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -9.45%  '

# Row 28
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.156'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -9.24%  '

# Row 29
$ws.Range('E29').Value = '  -10.76%  '

# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '7.12'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -4.55%  '

# Row 31
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.999'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -0.06%  '

# Row 32
$ws.Range('D32').Value = '0.0₃0714'
$ws.Range('E32').Value = '  -13.33%  '

# Row 33
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '146.49'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -3.12%  '

# Row 34
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '17.76'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -7.50%  '

# Row 35
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.44'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -9.99%  '

# Row 36
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '5.02'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -7.35%  '

# Row 37
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '3.55'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -14.91%  '

# Row 38
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.06'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -6.24%  '

# Row 39
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.800'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -15.11%  '

# Row 40
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.993'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -0.23%  '

# Row 41
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '32.98'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -8.29%  '

# Row 42
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.594'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -1.08%  '

# Row 43
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.0523'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -6.20%  '

# Row 44
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '3.26'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -8.20%  '

# Row 45
$ws.Range('B45').Value = 'Stacks'
$ws.Range('C45').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.24'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -10.60%  '

# Row 46
$ws.Range('B46').Value = 'WhiteBITCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '10.06'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -2.98%  '

# Row 47
$ws.Range('D47').Value = '1.933.48'
$ws.Range('E47').Value = '  -11.59%  '

# Row 48
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0886'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +0.26%  '

# Row 49
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0219'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -3.14%  '

# Row 50
$ws.Range('B50').Value = 'Bittensor'
$ws.Range('C50').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '233.19'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +3.30%  '

# Row 51
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '16.62'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -12.16%  '

